$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("T55")

# Row 2 (Q0)
$ws.Range("B2").Value = 0.1022639543468913
$ws.Range("C2").Value = 0.1022643483565146
$ws.Range("D2").Value = 0.230073309893088
$ws.Range("E2").Value = 0.479659577088885
$ws.Range("F2").Value = 0.4796595348963852
$ws.Range("G2").Value = 22

# Row 3 (Q1)
$ws.Range("B3").Value = 0.3374033227550764
$ws.Range("C3").Value = 2.647337262386587
$ws.Range("D3").Value = 27.11600199653427
$ws.Range("E3").Value = 5.207302756373425
$ws.Range("F3").Value = 5.324684877440589
$ws.Range("G3").Value = 21

# Row 4 (Q2)
$ws.Range("B4").Value = -1.028804989961556
$ws.Range("C4").Value = 2.587274768884802
$ws.Range("D4").Value = 20.70914113668103
$ws.Range("E4").Value = 4.550729736721467
$ws.Range("F4").Value = 4.548070895770996
$ws.Range("G4").Value = 20

# Row 5 (Q3)
$ws.Range("B5").Value = -0.2137546113718506
$ws.Range("C5").Value = 1.449643478764507
$ws.Range("D5").Value = 6.637417689940448
$ws.Range("E5").Value = 2.576318631291644
$ws.Range("F5").Value = 2.637789546666938
$ws.Range("G5").Value = 19

# Row 6 (Q4)
$ws.Range("B6").Value = -0.08368731705856663
$ws.Range("C6").Value = 1.582719152019898
$ws.Range("D6").Value = 10.41298270791714
$ws.Range("E6").Value = 3.22691535493529
$ws.Range("F6").Value = 3.319351677802833
$ws.Range("G6").Value = 18

# Row 7 (Q5)
$ws.Range("B7").Value = -0.40981972250463
$ws.Range("C7").Value = 2.002722685459293
$ws.Range("D7").Value = 13.90742892030138
$ws.Range("E7").Value = 3.729266539187214
$ws.Range("F7").Value = 3.820758303014839
$ws.Range("G7").Value = 17

# Row 8 (Q6)
$ws.Range("B8").Value = -0.2294455457529222
$ws.Range("C8").Value = 2.061615647959201
$ws.Range("D8").Value = 12.33890367925494
$ws.Range("E8").Value = 3.512677565512517
$ws.Range("F8").Value = 3.628201800890946
$ws.Range("G8").Value = 15
